$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FRAMECALCULATOR")
$ws.Range("F2").Value = 765
$ws.Range("F4").Value = 790
$ws.Range("K52").Value = 105
